# paises.xlsx data refresh: update country case counts as of 8 Apr 2020 18:22
# (replaces the prior 17:52 snapshot). Several countries climb in the
# "Casos totales" ranking because of the new numbers, which shuffles the
# sorted row order for: Arabia Saudita/Luxemburgo/Indonesia/Peru (rows 38-41),
# Cuba/Uruguay/Afganistan/Oman/Albania (rows 88-92),
# Gibraltar/Paraguay/Camboya (rows 125-127),
# Islas Caimanes/Congo (rows 144-145) and
# Malaui/Belice/Somalia (rows 193-195).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 8 de Abril de 2020 a las 18:22"

# Row 4: Estados Unidos (updated figures)
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 406585
$ws.Range("C4").Value = 6250
$ws.Range("D4").Value = 21993
$ws.Range("E4").Value = 371505
$ws.Range("F4").Value = 9220
$ws.Range("G4").Value = 246
$ws.Range("H4").Value = 13087

# Row 5: España (updated figures)
$ws.Range("A5").Value = "España"
$ws.Range("B5").Value = 146690
$ws.Range("C5").Value = 4748
$ws.Range("D5").Value = 48021
$ws.Range("E5").Value = 83996
$ws.Range("F5").Value = 7069
$ws.Range("G5").Value = 628
$ws.Range("H5").Value = 14673

# Row 6: Italia (updated figures)
$ws.Range("A6").Value = "Italia"
$ws.Range("B6").Value = 139422
$ws.Range("C6").Value = 3836
$ws.Range("D6").Value = 26491
$ws.Range("E6").Value = 95262
$ws.Range("F6").Value = 3693
$ws.Range("G6").Value = 542
$ws.Range("H6").Value = 17669

# Row 11: Reino Unido (updated figures)
$ws.Range("A11").Value = "Reino Unido"
$ws.Range("B11").Value = 60733
$ws.Range("C11").Value = 5491
$ws.Range("D11").Value = 135
$ws.Range("E11").Value = 53501
$ws.Range("F11").Value = 1559
$ws.Range("G11").Value = 938
$ws.Range("H11").Value = 7097

# Row 14: Suiza (updated figures)
$ws.Range("A14").Value = "Suiza"
$ws.Range("B14").Value = 23248
$ws.Range("C14").Value = 995
$ws.Range("D14").Value = 9800
$ws.Range("E14").Value = 12555
$ws.Range("F14").Value = 391
$ws.Range("G14").Value = 72
$ws.Range("H14").Value = 893

# Row 19: Austria (updated figures)
$ws.Range("A19").Value = "Austria"
$ws.Range("B19").Value = 12916
$ws.Range("C19").Value = 277
$ws.Range("D19").Value = 4512
$ws.Range("E19").Value = 8131
$ws.Range("F19").Value = 267
$ws.Range("G19").Value = 30
$ws.Range("H19").Value = 273

# Row 30: Chequia (updated figures)
$ws.Range("A30").Value = "Chequia"
$ws.Range("B30").Value = 5221
$ws.Range("C30").Value = 204
$ws.Range("D30").Value = 233
$ws.Range("E30").Value = 4889
$ws.Range("F30").Value = 103
$ws.Range("G30").Value = 11
$ws.Range("H30").Value = 99

# Row 38: Luxemburgo -> Arabia Saudita
$ws.Range("A38").Value = "Arabia Saudita"
$ws.Range("B38").Value = 3122
$ws.Range("C38").Value = 327
$ws.Range("D38").Value = 631
$ws.Range("E38").Value = 2450
$ws.Range("F38").Value = 41
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 41

# Row 39: Indonesia -> Luxemburgo
$ws.Range("A39").Value = "Luxemburgo"
$ws.Range("B39").Value = 3034
$ws.Range("C39").Value = 64
$ws.Range("D39").Value = 500
$ws.Range("E39").Value = 2488
$ws.Range("F39").Value = 34
$ws.Range("G39").Value = 2
$ws.Range("H39").Value = 46

# Row 40: Peru -> Indonesia
$ws.Range("A40").Value = "Indonesia"
$ws.Range("B40").Value = 2956
$ws.Range("C40").Value = 218
$ws.Range("D40").Value = 222
$ws.Range("E40").Value = 2494
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 19
$ws.Range("H40").Value = 240

# Row 41: Arabia Saudita -> Peru
$ws.Range("A41").Value = "Peru"
$ws.Range("B41").Value = 2954
$ws.Range("C41").Value = 0
$ws.Range("D41").Value = 1301
$ws.Range("E41").Value = 1546
$ws.Range("F41").Value = 109
$ws.Range("G41").Value = 0
$ws.Range("H41").Value = 107

# Row 57: Argelia (updated figures)
$ws.Range("A57").Value = "Argelia"
$ws.Range("B57").Value = 1572
$ws.Range("C57").Value = 104
$ws.Range("D57").Value = 237
$ws.Range("E57").Value = 1130
$ws.Range("F57").Value = 46
$ws.Range("G57").Value = 12
$ws.Range("H57").Value = 205

# Row 78: Eslovaquia (updated figures)
$ws.Range("A78").Value = "Eslovaquia"
$ws.Range("B78").Value = 682
$ws.Range("C78").Value = 101
$ws.Range("D78").Value = 16
$ws.Range("E78").Value = 664
$ws.Range("F78").Value = 3
$ws.Range("G78").Value = 0
$ws.Range("H78").Value = 2

# Row 88: Uruguay -> Cuba
$ws.Range("A88").Value = "Cuba"
$ws.Range("B88").Value = 457
$ws.Range("C88").Value = 61
$ws.Range("D88").Value = 27
$ws.Range("E88").Value = 418
$ws.Range("F88").Value = 15
$ws.Range("G88").Value = 1
$ws.Range("H88").Value = 12

# Row 89: Afganistan -> Uruguay
$ws.Range("A89").Value = "Uruguay"
$ws.Range("B89").Value = 424
$ws.Range("C89").Value = 0
$ws.Range("D89").Value = 150
$ws.Range("E89").Value = 267
$ws.Range("F89").Value = 14
$ws.Range("G89").Value = 0
$ws.Range("H89").Value = 7

# Row 90: Oman -> Afganistan
$ws.Range("A90").Value = "Afganistan"
$ws.Range("B90").Value = 423
$ws.Range("C90").Value = 0
$ws.Range("D90").Value = 18
$ws.Range("E90").Value = 391
$ws.Range("F90").Value = 0
$ws.Range("G90").Value = 0
$ws.Range("H90").Value = 14

# Row 91: Albania -> Oman
$ws.Range("A91").Value = "Oman"
$ws.Range("B91").Value = 419
$ws.Range("C91").Value = 48
$ws.Range("D91").Value = 72
$ws.Range("E91").Value = 345
$ws.Range("F91").Value = 3
$ws.Range("G91").Value = 0
$ws.Range("H91").Value = 2

# Row 92: Cuba -> Albania
$ws.Range("A92").Value = "Albania"
$ws.Range("B92").Value = 400
$ws.Range("C92").Value = 17
$ws.Range("D92").Value = 154
$ws.Range("E92").Value = 224
$ws.Range("F92").Value = 7
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 22

# Row 98: Ghana (updated figures)
$ws.Range("A98").Value = "Ghana"
$ws.Range("B98").Value = 313
$ws.Range("C98").Value = 26
$ws.Range("D98").Value = 34
$ws.Range("E98").Value = 273
$ws.Range("F98").Value = 2
$ws.Range("G98").Value = 1
$ws.Range("H98").Value = 6

# Row 125: Paraguay -> Gibraltar
$ws.Range("A125").Value = "Gibraltar"
$ws.Range("B125").Value = 120
$ws.Range("C125").Value = 7
$ws.Range("D125").Value = 60
$ws.Range("E125").Value = 60
$ws.Range("F125").Value = 0
$ws.Range("G125").Value = 0
$ws.Range("H125").Value = 0

# Row 126: Camboya -> Paraguay
$ws.Range("A126").Value = "Paraguay"
$ws.Range("B126").Value = 119
$ws.Range("C126").Value = 4
$ws.Range("D126").Value = 15
$ws.Range("E126").Value = 99
$ws.Range("F126").Value = 1
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = 5

# Row 127: Gibraltar -> Camboya
$ws.Range("A127").Value = "Camboya"
$ws.Range("B127").Value = 117
$ws.Range("C127").Value = 2
$ws.Range("D127").Value = 63
$ws.Range("E127").Value = 54
$ws.Range("F127").Value = 1
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 0

# Row 144: Congo -> Islas Caimanes
$ws.Range("A144").Value = "Islas Caimanes"
$ws.Range("B144").Value = 45
$ws.Range("C144").Value = 0
$ws.Range("D144").Value = 6
$ws.Range("E144").Value = 38
$ws.Range("F144").Value = 0
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 1

# Row 145: Islas Caimanes -> Congo
$ws.Range("A145").Value = "Congo"
$ws.Range("B145").Value = 45
$ws.Range("C145").Value = 0
$ws.Range("D145").Value = 2
$ws.Range("E145").Value = 38
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 5

# Row 193: Somalia -> Malaui
$ws.Range("A193").Value = "Malaui"
$ws.Range("B193").Value = 8
$ws.Range("C193").Value = 0
$ws.Range("D193").Value = 0
$ws.Range("E193").Value = 7
$ws.Range("F193").Value = 1
$ws.Range("G193").Value = 0
$ws.Range("H193").Value = 1

# Row 194: Malaui -> Belice
$ws.Range("A194").Value = "Belice"
$ws.Range("B194").Value = 8
$ws.Range("C194").Value = 1
$ws.Range("D194").Value = 0
$ws.Range("E194").Value = 7
$ws.Range("F194").Value = 1
$ws.Range("G194").Value = 0
$ws.Range("H194").Value = 1

# Row 195: Belice -> Somalia
$ws.Range("A195").Value = "Somalia"
$ws.Range("B195").Value = 8
$ws.Range("C195").Value = 0
$ws.Range("D195").Value = 1
$ws.Range("E195").Value = 6
$ws.Range("F195").Value = 0
$ws.Range("G195").Value = 1
$ws.Range("H195").Value = 1

